$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.59"
$ws.Range("E2").Value = "'0.60%"
$ws.Range("D3").Value = "'29.91"
$ws.Range("E3").Value = "'12.51%"
$ws.Range("D4").Value = "'5.144"
$ws.Range("E4").Value = "'0.18%"
$ws.Range("D5").Value = "'0.05732"
$ws.Range("E5").Value = "'2.16%"
$ws.Range("D6").Value = "'6.576"
$ws.Range("E6").Value = "'1.72%"
$ws.Range("D7").Value = "'0.8570"
$ws.Range("E7").Value = "'4.63%"
$ws.Range("D8").Value = "'0.8692"
$ws.Range("E8").Value = "'4.31%"
$ws.Range("D9").Value = "'0.1354"
$ws.Range("E9").Value = "'1.90%"
$ws.Range("D10").Value = "'0.06925"
$ws.Range("E10").Value = "'-0.01%"
$ws.Range("D11").Value = "'0.02905"
$ws.Range("E11").Value = "'0.35%"
$ws.Range("D12").Value = "'0.09363"
$ws.Range("E12").Value = "'-0.23%"
$ws.Range("D13").Value = "'0.001518"
$ws.Range("E13").Value = "'-0.11%"
$ws.Range("D14").Value = "'0.04152"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005986"
$ws.Range("E15").Value = "'0.12%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006063"
$ws.Range("E16").Value = "'-1.45%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.509"
$ws.Range("E17").Value = "'-3.79%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.018"
$ws.Range("E18").Value = "'-0.07%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.175"
$ws.Range("E19").Value = "'-0.70%"
$ws.Range("D20").Value = "'0.3142"
$ws.Range("E20").Value = "'0.98%"
$ws.Range("D21").Value = "'0.03308"
$ws.Range("E21").Value = "'8.01%"
$ws.Range("D22").Value = "'0.1304"
$ws.Range("E22").Value = "'0.35%"
$ws.Range("D23").Value = "'3.602"
$ws.Range("E23").Value = "'-4.18%"
$ws.Range("E24").Value = "'2.56%"
$ws.Range("D25").Value = "'0.001209"
$ws.Range("E25").Value = "'-1.39%"
$ws.Range("D26").Value = "'0.004482"
$ws.Range("E27").Value = "'22.68%"
$ws.Range("E28").Value = "'-0.78%"
$ws.Range("D40").Value = "'0.03772"
$ws.Range("E40").Value = "'3.66%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1066"
$ws.Range("E41").Value = "'-22.24%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.003498"
$ws.Range("E42").Value = "'-43.30%"
$ws.Range("D43").Value = "'0.002366"
$ws.Range("E43").Value = "'-8.64%"
$ws.Range("D44").Value = "'0.009877"
$ws.Range("E44").Value = "'21.62%"
$ws.Range("D45").Value = "'0.00005080"
$ws.Range("E45").Value = "'-5.02%"
$ws.Range("E46").Value = "'-0.22%"
$ws.Range("D47").Value = "'0.07982"
$ws.Range("E47").Value = "'-26.77%"
$ws.Range("D48").Value = "'0.002737"
$ws.Range("E48").Value = "'8.40%"
$ws.Range("E49").Value = "'-0.22%"
$ws.Range("E50").Value = "'-0.22%"
